# Automatische test-sync: 2025-08-18 21:25:50
# Appends the newest mail-log entry to the "Logs" sheet and refreshes the
# category counts on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new log row (row 11) on the "Logs" sheet -----------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A11").Value = "Geen onderwerp"
$logs.Range("B11").Value = "onbekend"
$logs.Range("D11").Value = "Overig"
$logs.Range("F11").Value = "2025-08-18 21:25:44"
$logs.Range("G11").Value = "Nee"
$logs.Range("H11").Value = "Ja"
$logs.Range("I11").Value = "Nee"
$logs.Range("J11").Value = "Nee"

# --- 2. Extend the conditional formatting ranges to include the new row -----
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D11"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G11"))
$logs.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H11"))
$logs.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I11"))
$logs.Range("J2:J10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J11"))

# --- 3. Refresh the "Dashboard" category counts ------------------------------
# The new row's category ("Overig") now has 3 entries and overtakes
# "Intern verzoek / Actie voor medewerker" (2 entries), so the two rows swap.
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Overig"
$dash.Range("B3").Value = 3
$dash.Range("A4").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B4").Value = 2
